$wb = $excel.ActiveWorkbook

$ws = $wb.Sheets.Item("ALC")
$ws.Range("H19").Value = 4747.44
$ws.Range("I19").Value = 9361.546
$ws.Range("J19").Value = 1122.0714
$ws.Range("K19").Value = 9361.546
$ws.Range("L19").Value = 1122.0714
$ws.Range("M19").Value = -9186.546
$ws.Range("N19").Value = -1472.0714
$ws.Range("H106").Value = 3956.6667
$ws.Range("I106").Value = 1758.3334
$ws.Range("J106").Value = 5422.222
$ws.Range("K106").Value = 1758.3334
$ws.Range("L106").Value = 5422.222
$ws.Range("M106").Value = -1127.3334
$ws.Range("N106").Value = -6684.222
$ws.Range("H116").Value = 5733.3335
$ws.Range("I116").Value = 6175
$ws.Range("K116").Value = 6175
$ws.Range("M116").Value = -2733
$ws.Range("H132").Value = 2792.9412
$ws.Range("I132").Value = 1698.575
$ws.Range("K132").Value = 5095.725
$ws.Range("M132").Value = -2565.725
$ws.Range("H137").Value = 2814
$ws.Range("I137").Value = 2760.1924
$ws.Range("J137").Value = 2941.182
$ws.Range("K137").Value = 8280.5772
$ws.Range("L137").Value = 8823.545999999998
$ws.Range("M137").Value = -5730.5772
$ws.Range("N137").Value = -13923.546
$ws.Range("H138").Value = 2201.1409
$ws.Range("I138").Value = 1026.4517
$ws.Range("J138").Value = 3111.525
$ws.Range("K138").Value = 3079.3551
$ws.Range("L138").Value = 9334.575000000001
$ws.Range("M138").Value = 2060.6449
$ws.Range("N138").Value = -19614.575

$ws = $wb.Sheets.Item("ARM")
$ws.Range("H32").Value = 1348084.2
$ws.Range("I32").Value = 1576275.5
$ws.Range("J32").Value = 4291.5557
$ws.Range("K32").Value = 1576275.5
$ws.Range("L32").Value = 4291.5557
$ws.Range("M32").Value = -1575988.5
$ws.Range("N32").Value = -4865.5557
$ws.Range("H61").Value = 16727.688
$ws.Range("I61").Value = 22095.396
$ws.Range("J61").Value = 3167.158
$ws.Range("K61").Value = 22095.396
$ws.Range("L61").Value = 3167.158
$ws.Range("M61").Value = -21883.396
$ws.Range("N61").Value = -3591.158
$ws.Range("H110").Value = 934.4167
$ws.Range("I110").Value = 914.125
$ws.Range("J110").Value = 975
$ws.Range("K110").Value = 914.125
$ws.Range("L110").Value = 975
$ws.Range("M110").Value = 1130.875
$ws.Range("N110").Value = -5065
$ws.Range("H122").Value = 3039.457
$ws.Range("I122").Value = 2805.5908
$ws.Range("J122").Value = 3435.2307
$ws.Range("K122").Value = 8416.7724
$ws.Range("L122").Value = 10305.6921
$ws.Range("M122").Value = -5966.7724
$ws.Range("N122").Value = -15205.6921
$ws.Range("H132").Value = 26660.488
$ws.Range("I132").Value = 49047.637
$ws.Range("J132").Value = 3207.2856
$ws.Range("K132").Value = 147142.911
$ws.Range("L132").Value = 9621.856800000001
$ws.Range("M132").Value = -144612.911
$ws.Range("N132").Value = -14681.8568
$ws.Range("H136").Value = 16727.688
$ws.Range("I136").Value = 22095.396
$ws.Range("J136").Value = 3167.158
$ws.Range("K136").Value = 66286.18799999999
$ws.Range("L136").Value = 9501.474
$ws.Range("M136").Value = -63736.18799999999
$ws.Range("N136").Value = -14601.474

$ws = $wb.Sheets.Item("BSM")
$ws.Range("H80").Value = 141
$ws.Range("I80").Value = 161.5
$ws.Range("J80").Value = 100
$ws.Range("K80").Value = 161.5
$ws.Range("L80").Value = 100
$ws.Range("M80").Value = 836.5
$ws.Range("N80").Value = -2096
$ws.Range("H83").Value = 141
$ws.Range("I83").Value = 161.5
$ws.Range("J83").Value = 100
$ws.Range("K83").Value = 807.5
$ws.Range("L83").Value = 500
$ws.Range("M83").Value = 4184.5
$ws.Range("N83").Value = -10484
$ws.Range("H105").Value = 1497.7391
$ws.Range("I105").Value = 1449.9048
$ws.Range("J105").Value = 2000
$ws.Range("K105").Value = 1449.9048
$ws.Range("L105").Value = 2000
$ws.Range("M105").Value = 297.0952
$ws.Range("N105").Value = -5494
$ws.Range("H134").Value = 2593.7322
$ws.Range("I134").Value = 2340.5813
$ws.Range("J134").Value = 3431.077
$ws.Range("K134").Value = 7021.743899999999
$ws.Range("L134").Value = 10293.231
$ws.Range("M134").Value = -4486.743899999999
$ws.Range("N134").Value = -15363.231

$ws = $wb.Sheets.Item("CRP")
$ws.Range("H31").Value = 2749.3103
$ws.Range("I31").Value = 1971.75
$ws.Range("J31").Value = 4477.222
$ws.Range("K31").Value = 1971.75
$ws.Range("L31").Value = 4477.222
$ws.Range("M31").Value = -1676.75
$ws.Range("N31").Value = -5067.222
$ws.Range("H34").Value = 2749.3103
$ws.Range("I34").Value = 1971.75
$ws.Range("J34").Value = 4477.222
$ws.Range("K34").Value = 1971.75
$ws.Range("L34").Value = 4477.222
$ws.Range("M34").Value = -1769.75
$ws.Range("N34").Value = -4881.222
$ws.Range("H58").Value = 1573.5
$ws.Range("I58").Value = 804.5454999999999
$ws.Range("J58").Value = 2983.25
$ws.Range("K58").Value = 804.5454999999999
$ws.Range("L58").Value = 2983.25
$ws.Range("M58").Value = -601.5454999999999
$ws.Range("N58").Value = -3389.25
$ws.Range("H99").Value = 65030.25
$ws.Range("I99").Value = 101696.9
$ws.Range("K99").Value = 101696.9
$ws.Range("M99").Value = -100198.9
$ws.Range("H105").Value = 1103.6666
$ws.Range("I105").Value = 705.7143
$ws.Range("K105").Value = 705.7143
$ws.Range("M105").Value = 1041.2857
$ws.Range("H126").Value = 65030.25
$ws.Range("I126").Value = 101696.9
$ws.Range("K126").Value = 305090.7
$ws.Range("M126").Value = -302620.7
$ws.Range("H132").Value = 2608.9355
$ws.Range("I132").Value = 1321.6
$ws.Range("J132").Value = 3815.8125
$ws.Range("K132").Value = 3964.8
$ws.Range("L132").Value = 11447.4375
$ws.Range("M132").Value = -1434.8
$ws.Range("N132").Value = -16507.4375
$ws.Range("H136").Value = 1573.5
$ws.Range("I136").Value = 804.5454999999999
$ws.Range("J136").Value = 2983.25
$ws.Range("K136").Value = 2413.6365
$ws.Range("L136").Value = 8949.75
$ws.Range("M136").Value = 136.3635000000004
$ws.Range("N136").Value = -14049.75

$ws = $wb.Sheets.Item("CUL")
$ws.Range("H132").Value = 6668
$ws.Range("I132").Value = 4031.3333
$ws.Range("J132").Value = 8250
$ws.Range("K132").Value = 36281.9997
$ws.Range("L132").Value = 74250
$ws.Range("M132").Value = -33751.9997
$ws.Range("N132").Value = -79310
$ws.Range("H137").Value = 3100.0417
$ws.Range("I137").Value = 1292.8572
$ws.Range("J137").Value = 5630.1
$ws.Range("K137").Value = 3878.5716
$ws.Range("L137").Value = 16890.3
$ws.Range("M137").Value = 1221.4284
$ws.Range("N137").Value = -27090.3

$ws = $wb.Sheets.Item("GSM")
$ws.Range("H102").Value = 3380.0715
$ws.Range("I102").Value = 1324.3334
$ws.Range("K102").Value = 1324.3334
$ws.Range("M102").Value = 297.6666
$ws.Range("H122").Value = 1713.6111
$ws.Range("I122").Value = 1798.5555
$ws.Range("J122").Value = 1628.6666
$ws.Range("K122").Value = 5395.666499999999
$ws.Range("L122").Value = 4885.9998
$ws.Range("M122").Value = -2945.666499999999
$ws.Range("N122").Value = -9785.9998
$ws.Range("H132").Value = 3273.182
$ws.Range("I132").Value = 3107.7073
$ws.Range("J132").Value = 3757.7856
$ws.Range("K132").Value = 9323.1219
$ws.Range("L132").Value = 11273.3568
$ws.Range("M132").Value = -6793.1219
$ws.Range("N132").Value = -16333.3568

$ws = $wb.Sheets.Item("LTW")
$ws.Range("H82").Value = 2216.5
$ws.Range("I82").Value = 2114
$ws.Range("J82").Value = 2319
$ws.Range("K82").Value = 2114
$ws.Range("L82").Value = 2319
$ws.Range("M82").Value = -1753
$ws.Range("N82").Value = -3041
$ws.Range("H85").Value = 2216.5
$ws.Range("I85").Value = 2114
$ws.Range("J85").Value = 2319
$ws.Range("K85").Value = 2114
$ws.Range("L85").Value = 2319
$ws.Range("M85").Value = -866
$ws.Range("N85").Value = -4815

$ws = $wb.Sheets.Item("WVR")
$ws.Range("H107").Value = 608.8182
$ws.Range("I107").Value = 596.63635
$ws.Range("J107").Value = 621
$ws.Range("K107").Value = 1789.90905
$ws.Range("L107").Value = 1863
$ws.Range("M107").Value = 130.09095
$ws.Range("N107").Value = -5703
$ws.Range("H110").Value = 44000
$ws.Range("J110").Value = 44000
$ws.Range("L110").Value = 44000
$ws.Range("N110").Value = -52180
$ws.Range("H118").Value = 27950.5
$ws.Range("J118").Value = 27950.5
$ws.Range("L118").Value = 27950.5
$ws.Range("N118").Value = -31264.5
$ws.Range("H121").Value = 30210
$ws.Range("J121").Value = 30210
$ws.Range("L121").Value = 30210
$ws.Range("N121").Value = -33704
$ws.Range("H135").Value = 53000
$ws.Range("J135").Value = 53000
$ws.Range("L135").Value = 53000
$ws.Range("N135").Value = -63140
$ws.Range("H137").Value = 47921.25
$ws.Range("I137").Value = 0
$ws.Range("J137").Value = 47921.25
$ws.Range("K137").Value = 0
$ws.Range("L137").Value = 47921.25
$ws.Range("M137").ClearContents()
$ws.Range("N137").Value = -58121.25
